$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 149 values ---
$ws.Range("D149").Value = 0.1694277792
$ws.Range("E149").Value = 0.1920179357
$ws.Range("F149").Value = -0.002234255455634382
$ws.Range("G149").Value = 0.1281540008926507

# --- Add new rows 150 and 151, copying formatting from row 149 ---
$ws.Range("A149:G149").Copy()
$ws.Range("A150:G151").PasteSpecial(-4122)

# Row 150
$ws.Range("A150").Value = 148
$ws.Range("B150").Value = 148
$ws.Range("C150").Value = 44557
$ws.Range("D150").Value = 0.1899679581
$ws.Range("E150").Value = 0.1745944373
$ws.Range("F150").Value = 0.1212326514399595
$ws.Range("G150").Value = -0.0907389111151663

# Row 151
$ws.Range("A151").Value = 149
$ws.Range("B151").Value = 149
$ws.Range("C151").Value = 44564
$ws.Range("D151").Value = 0.1745684611
$ws.Range("E151").Value = 0
$ws.Range("F151").Value = -0.08106365491328615
$ws.Range("G151").Value = -1
